# Atualização de bases das ligas, do dia: 30-05-2024 às 12:21
#
# Refreshes the 5 match records (rows 309-313) on the "Romania Liga I"
# sheet with corrected data from the upstream league database re-sync:
# match id (B), home/away teams (E/F), full/half-time score (G:J),
# result (K) and the odds/closing-line columns (L:AD). Only the cells
# whose values actually changed are touched; Div (C) and Date (D) are
# unchanged for all five rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(309, 2).Value2 = 8191476
$ws.Cells.Item(309, 5).Value2 = "FC Voluntari"
$ws.Cells.Item(309, 6).Value2 = "Universitatea Cluj"
$ws.Cells.Item(309, 7).Value2 = 0
$ws.Cells.Item(309, 8).Value2 = 1
$ws.Cells.Item(309, 10).Value2 = 1
$ws.Cells.Item(309, 11).Value2 = "A"
$ws.Cells.Item(309, 12).Value2 = 3.05
$ws.Cells.Item(309, 14).Value2 = 2.15
$ws.Cells.Item(309, 15).Value2 = 2.6
$ws.Cells.Item(309, 16).Value2 = 3.4
$ws.Cells.Item(309, 17).Value2 = 2.4
$ws.Cells.Item(309, 18).Value2 = 0
$ws.Cells.Item(309, 19).Value2 = 2
$ws.Cells.Item(309, 20).Value2 = 1.85
$ws.Cells.Item(309, 22).Value2 = 2
$ws.Cells.Item(309, 23).Value2 = 1.85
$ws.Cells.Item(309, 24).Value2 = -1
$ws.Cells.Item(309, 26).Value2 = 1.4
$ws.Cells.Item(309, 27).Value2 = -1
$ws.Cells.Item(309, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(309, 29).Value2 = -1
$ws.Cells.Item(309, 30).Value2 = 0.8500000000000001
$ws.Cells.Item(310, 2).Value2 = 8191523
$ws.Cells.Item(310, 5).Value2 = "Otelul Galati"
$ws.Cells.Item(310, 6).Value2 = "FC Botosani"
$ws.Cells.Item(310, 7).Value2 = 2
$ws.Cells.Item(310, 8).Value2 = 0
$ws.Cells.Item(310, 9).Value2 = 2
$ws.Cells.Item(310, 10).Value2 = 0
$ws.Cells.Item(310, 11).Value2 = "H"
$ws.Cells.Item(310, 12).Value2 = 1.666
$ws.Cells.Item(310, 13).Value2 = 3.6
$ws.Cells.Item(310, 14).Value2 = 4.6
$ws.Cells.Item(310, 15).Value2 = 2.9
$ws.Cells.Item(310, 16).Value2 = 3.5
$ws.Cells.Item(310, 17).Value2 = 2.2
$ws.Cells.Item(310, 18).Value2 = 0.25
$ws.Cells.Item(310, 19).Value2 = 1.85
$ws.Cells.Item(310, 20).Value2 = 2
$ws.Cells.Item(310, 22).Value2 = 1.875
$ws.Cells.Item(310, 23).Value2 = 1.975
$ws.Cells.Item(310, 24).Value2 = 1.9
$ws.Cells.Item(310, 26).Value2 = -1
$ws.Cells.Item(310, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(310, 28).Value2 = -1
$ws.Cells.Item(310, 29).Value2 = -0.5
$ws.Cells.Item(310, 30).Value2 = 0.4875
$ws.Cells.Item(311, 2).Value2 = 8191462
$ws.Cells.Item(311, 5).Value2 = "CSM Politehnica Iasi"
$ws.Cells.Item(311, 6).Value2 = "Petrolul Ploiesti"
$ws.Cells.Item(311, 7).Value2 = 2
$ws.Cells.Item(311, 8).Value2 = 0
$ws.Cells.Item(311, 11).Value2 = "H"
$ws.Cells.Item(311, 12).Value2 = 2.1
$ws.Cells.Item(311, 14).Value2 = 3.1
$ws.Cells.Item(311, 15).Value2 = 1.8
$ws.Cells.Item(311, 16).Value2 = 3.2
$ws.Cells.Item(311, 17).Value2 = 4.2
$ws.Cells.Item(311, 18).Value2 = -0.5
$ws.Cells.Item(311, 22).Value2 = 2.025
$ws.Cells.Item(311, 23).Value2 = 1.825
$ws.Cells.Item(311, 24).Value2 = 0.8
$ws.Cells.Item(311, 26).Value2 = -1
$ws.Cells.Item(311, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(311, 28).Value2 = -1
$ws.Cells.Item(311, 29).Value2 = -0.5
$ws.Cells.Item(311, 30).Value2 = 0.4125
$ws.Cells.Item(312, 2).Value2 = 8191463
$ws.Cells.Item(312, 5).Value2 = "Dinamo Bucharest"
$ws.Cells.Item(312, 6).Value2 = "ACS UTA Batrana Doamna"
$ws.Cells.Item(312, 12).Value2 = 1.833
$ws.Cells.Item(312, 13).Value2 = 3.4
$ws.Cells.Item(312, 14).Value2 = 3.6
$ws.Cells.Item(312, 15).Value2 = 1.5
$ws.Cells.Item(312, 16).Value2 = 4.333
$ws.Cells.Item(312, 17).Value2 = 5
$ws.Cells.Item(312, 18).Value2 = -1
$ws.Cells.Item(312, 19).Value2 = 1.875
$ws.Cells.Item(312, 20).Value2 = 1.975
$ws.Cells.Item(312, 21).Value2 = 3
$ws.Cells.Item(312, 22).Value2 = 2.025
$ws.Cells.Item(312, 23).Value2 = 1.825
$ws.Cells.Item(312, 24).Value2 = 0.5
$ws.Cells.Item(312, 27).Value2 = 0.875
$ws.Cells.Item(312, 29).Value2 = -1
$ws.Cells.Item(312, 30).Value2 = 0.825
$ws.Cells.Item(313, 2).Value2 = 8191475
$ws.Cells.Item(313, 5).Value2 = "FC U Craiova 1948"
$ws.Cells.Item(313, 6).Value2 = "AFC Hermannstadt"
$ws.Cells.Item(313, 7).Value2 = 1
$ws.Cells.Item(313, 8).Value2 = 3
$ws.Cells.Item(313, 9).Value2 = 0
$ws.Cells.Item(313, 11).Value2 = "A"
$ws.Cells.Item(313, 12).Value2 = 2.625
$ws.Cells.Item(313, 13).Value2 = 3.3
$ws.Cells.Item(313, 14).Value2 = 2.45
$ws.Cells.Item(313, 15).Value2 = 2.05
$ws.Cells.Item(313, 16).Value2 = 3.5
$ws.Cells.Item(313, 17).Value2 = 3
$ws.Cells.Item(313, 18).Value2 = -0.25
$ws.Cells.Item(313, 19).Value2 = 1.85
$ws.Cells.Item(313, 20).Value2 = 2
$ws.Cells.Item(313, 21).Value2 = 2.25
$ws.Cells.Item(313, 22).Value2 = 1.825
$ws.Cells.Item(313, 23).Value2 = 2.025
$ws.Cells.Item(313, 24).Value2 = -1
$ws.Cells.Item(313, 26).Value2 = 2
$ws.Cells.Item(313, 27).Value2 = -1
$ws.Cells.Item(313, 28).Value2 = 1
$ws.Cells.Item(313, 29).Value2 = 0.825
$ws.Cells.Item(313, 30).Value2 = -1
